# "new pagination from DB directly"
# Add an "issue" tracking column (I) to the feature sheet, move the
# existing long-form issue note from E17 into I17, mark the affected rows
# as "suspend", mark the pagination rows (6 & 26) as "done" with the new
# technique note, and record a new failure note on row 37.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("feature")

# --- New header for the issue column ---
$ws.Range("I1").Value = "issue"

# --- Row 6: pagination work is now done, drop the old highlight, note the tech used ---
$ws.Range("E7").Copy()
$ws.Range("E6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E6").Value = "done"

$ws.Range("F7").Copy()
$ws.Range("F6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F6").Value = 43986

$ws.Range("G6").Value = "collections.sort or JDBC order by"

# --- Row 17: suspend the order placement feature, keep the detailed issue note but move it to column I ---
$ws.Range("E17").Value = "suspend"

$ws.Range("I17").Value = "issue ,can not populate order data into db"
$ws.Range("E17").Copy()
$ws.Range("I17").PasteSpecial(-4122)  # xlPasteFormats (reapply after value overwrite below)
$ws.Range("I17").Value = "issue ,can not populate order data into db"

# --- Row 26: search-book-by-criteria pagination is done too ---
$ws.Range("E26").Value = "done"
$ws.Range("F26").Value = 43986
$ws.Range("G26").Value = "collections.sort or JDBC order by"

# --- Row 37: importing js currently fails, suspend and record the error ---
$ws.Range("E18").Copy()
$ws.Range("E37").PasteSpecial(-4122)  # xlPasteFormats (same yellow highlight style as other "suspend"/issue cells)
$ws.Range("E37").Value = "suspend"

$ws.Range("I37").Value = "cannot do it , error"

$ws.Range("E9").Select()
